$wb = $excel.ActiveWorkbook

# Rename Sheet2 -> InvalidLogin
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "InvalidLogin"

# Populate the InvalidLogin sheet with data-driven test data
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "admin"

# Update selection on the first sheet (ValidLogin) to A1:B2
$ws1 = $wb.Worksheets.Item("ValidLogin")
$ws1.Range("A1:B2").Select()

# Make InvalidLogin the active sheet/tab, with C4 selected
$ws2.Activate()
$ws2.Range("C4").Select()
